# Add new materials to "Materiais" sheet and new recipe lines to "Receitas" sheet,
# matching the data uploaded in the author's spreadsheet update.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Materiais")
$ws2 = $wb.Worksheets.Item("Receitas")

# ---- Materiais: new rows 3-11 ----
$materiais = @(
  ,@("Dióxido de Titânio", 18.2, 4.2300000000000004, "13463-67-7", "irritante")
  ,@("Acetato de Etila", 8.9, 0.9, "141-78-6", "inflamável")
  ,@("Xileno", 6.5, 0.86, "1330-20-7", "nocivo/inflamável")
  ,@("Carbonato de Cálcio", 2.1, 2.71, "471-34-1", "nenhum")
  ,@("Resina Epóxi", 22, 1.1599999999999999, "25068-38-6", "irritante/sensibilizante")
  ,@("Etanol 96%", 4.8, 0.81, "64-17-5", "inflamável")
  ,@("Pigmento Azul Ftalo", 45, 1.6, "147-14-8", "nenhum")
  ,@("Talco Industrial", 1.5, 2.75, "14807-96-6", "inalação perigosa")
  ,@("Aditivo Antiespumante", 32.5, 0.98, "63148-62-9", "nenhum")
)

$r = 3
foreach ($row in $materiais) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).NumberFormat = "@"
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---- Receitas: new rows 3-12 ----
$receitas = @(
  ,@("Tinta Azul Industrial", "Pigmento Azul Ftalo", 5)
  ,@("Tinta Azul Industrial", "Xileno", 30)
  ,@("Tinta Azul Industrial", "Aditivo Antiespumante", "0.5")
  ,@("Tinta Azul Industrial", "Carbonato de Cálcio", 10)
  ,@("Primer Branco Epóxi", "Resina Epóxi", 40)
  ,@("Primer Branco Epóxi", "Dióxido de Titânio", 15)
  ,@("Primer Branco Epóxi", "Acetato de Etila", 10)
  ,@("Primer Branco Epóxi", "Talco Industrial", 25)
  ,@("Solução de Limpeza", "Etanol 96%", 85)
  ,@("Solução de Limpeza", "Acetato de Etila", 15)
)

$r = 3
foreach ($row in $receitas) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $qty = $row[2]
    if ($qty -eq "0.5") {
        # This quantity was typed as text in the original sheet (right-aligned, General format)
        $ws2.Cells.Item($r, 3).Value = "'0.5"
        $ws2.Cells.Item($r, 3).HorizontalAlignment = -4152
    } else {
        $ws2.Cells.Item($r, 3).Value = $qty
    }
    $r = $r + 1
}

# ---- Sheet view / selection bookkeeping ----
# Receitas keeps its own pending selection, but Materiais becomes the active tab.
$ws2.Range("I17").Select()
$ws1.Activate()
$ws1.Range("D7").Select()

# ---- Print setup on Materiais ----
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1
